$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
$ws.Range("D2").Value = "64.957.55"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "3.512.04"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  -2.19%  "
$ws.Range("E6").Value = "  -1.25%  "
$ws.Range("D7").Value = "3.511.03"
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +3.08%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "4.111.31"
$ws.Range("E13").Value = "  -1.35%  "
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  +0.84%  "
$ws.Range("D17").Value = "3.512.99"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("D18").Value = "64.955.66"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("E23").Value = "  -0.41%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").Value = "3.654.39"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("E28").Value = "  +7.88%  "
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("E31").Value = "  -1.35%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "3.520.60"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("E37").Value = "  +4.36%  "
$ws.Range("E38").Value = "  +1.99%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("E43").Value = "  -2.73%  "
$ws.Range("E44").Value = "  +0.69%  "
$ws.Range("E45").Value = "  +3.32%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("D49").Value = "2.471.64"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E51").Value = "  +2.58%  "

# Numeric-looking values that must stay text: force Text format, assign, then
# clear the format again so the cell keeps its original (default) style.
$numericCells = @("D5","D6","D11","D14","D19","D20","D22","D23","D24","D27","D28","D29","D32","D34","D37","D38","D39","D40","D41","D42","D43","D44","D47","D48","D50","D51")
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}
$ws.Range("D5").Value = "585.80"
$ws.Range("D6").Value = "133.49"
$ws.Range("D11").Value = "7.16"
$ws.Range("D14").Value = "27.54"
$ws.Range("D19").Value = "10.01"
$ws.Range("D20").Value = "14.28"
$ws.Range("D22").Value = "390.39"
$ws.Range("D23").Value = "0.574"
$ws.Range("D24").Value = "74.79"
$ws.Range("D27").Value = "0.0000110"
$ws.Range("D28").Value = "1.60"
$ws.Range("D29").Value = "7.55"
$ws.Range("D32").Value = "8.26"
$ws.Range("D34").Value = "24.06"
$ws.Range("D37").Value = "5.17"
$ws.Range("D38").Value = "1.56"
$ws.Range("D39").Value = "169.42"
$ws.Range("D40").Value = "6.92"
$ws.Range("D41").Value = "0.0804"
$ws.Range("D42").Value = "0.818"
$ws.Range("D43").Value = "26.04"
$ws.Range("D44").Value = "42.94"
$ws.Range("D47").Value = "4.41"
$ws.Range("D48").Value = "1.64"
$ws.Range("D50").Value = "6.84"
$ws.Range("D51").Value = "0.0267"
foreach ($addr in $numericCells) {
    $ws.Range($addr).ClearFormats()
}
